$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("first_eval")

# New values per row (row -> B,C,D,E,F,G)
$data = @{
    2  = @(0.372432736606822,   2.619063551738779,  23.15927064259555,  4.812407987961489, 4.845717262623262, 51)
    3  = @(0.06436319397563839, 2.744095017829264,  23.38972601463277,  4.836292589849458, 4.884960613746139, 50)
    4  = @(0.3367729254375606,  2.905511120774942,  22.73626976693513,  4.768256470339566, 4.8056386898139,   49)
    5  = @(0.1606920638432899,  3.013415626553616,  24.31814105163908,  4.931342722995338, 4.98088114350354,  48)
    6  = @(0.3452801118255731,  3.332856983624331,  25.17271329514386,  5.017241602229642, 5.059460063471103, 47)
    7  = @(0.2526637188661144,  3.173335729886531,  26.35108212670401,  5.133330510176021, 5.183763510351113, 46)
    8  = @(0.3086035603155297,  3.48111105990218,   26.05187013802933,  5.104103264828145, 5.152335148274883, 45)
    9  = @(0.2033787227829306,  3.124439938622011,  23.6251089138144,   4.860566727637262, 4.912454093265776, 44)
    10 = @(0.2285605960393908,  3.349377590055617,  25.07901460685525,  5.007895227224233, 5.061882085544378, 43)
    11 = @(0.3835128682399718,  3.331810460747813,  25.45961486176948,  5.045752160161008, 5.092142141478524, 42)
}

$columns = @("B", "C", "D", "E", "F", "G")

foreach ($row in $data.Keys) {
    $values = $data[$row]
    for ($i = 0; $i -lt $columns.Length; $i++) {
        $col = $columns[$i]
        $ws.Range("$col$row").Value = $values[$i]
    }
}

$wb.Save()
